# Add "Garantías Propuestas" guarantee columns (AB:AF) to the export sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1) ---------------------------------------------
$ws.Range("AB1").Value = "Garantías Propuestas"
$ws.Range("AC1").Value = "Tipo Garantía"
$ws.Range("AD1").Value = "Tipo bien"
$ws.Range("AE1").Value = "Valor Comercial"
$ws.Range("AF1").Value = "Descripcion"

# Wrap the two longer headers so the taller header row reads cleanly.
$ws.Range("AB1").WrapText = $true
$ws.Range("AE1").WrapText = $true

# Header row is taller now that some cells wrap to two lines.
$ws.Rows.Item(1).RowHeight = 28.35

# --- New data row (row 2) -------------------------------------------------
$ws.Range("AB2").Value = $ws.Range("H2").Value2
$ws.Range("AC2").Value = "A13"
$ws.Range("AD2").Value = "110"
$ws.Range("AE2").Value = "120000"
$ws.Range("AF2").Value = "TERRENO"

# --- View state: scroll over to the new columns and select AF2 -----------
$excel.ActiveWindow.ScrollColumn = 27
$ws.Range("AF2").Select()
